$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the sample name / suffix text used throughout rows 2-5
# (these cells share strings with rows 6-9, which are being cleared below,
# so re-assign the same new text to each surviving cell)
$ws.Range("A2").Value = "cerebellum_tile1"
$ws.Range("A3").Value = "cerebellum_tile1"
$ws.Range("A4").Value = "cerebellum_tile1"
$ws.Range("A5").Value = "cerebellum_tile1"

$ws.Range("B2").Value = "039_GMB_tileRingMixScan_4rings_7scans"
$ws.Range("B3").Value = "039_GMB_tileRingMixScan_4rings_7scans"
$ws.Range("B4").Value = "039_GMB_tileRingMixScan_4rings_7scans"
$ws.Range("B5").Value = "039_GMB_tileRingMixScan_4rings_7scans"

# Clear out the data rows 6-9 entirely (only the styled, empty D cell remains)
$ws.Range("A6:E9").ClearContents()

# Update the active selection to A3:A5 (active cell A3)
$ws.Range("A3:A5").Select()
